# EPS v3.3.1 -> v3.4.2 update for "Hydgn Production Eff by Pathway.xlsx"
#
# The only substantive change in this revision is on the "HPEbP" sheet:
# the natural-gas-reforming efficiency formula in B3 drops the "+46"
# term (118/(162+2+46) -> 118/(162+2)); the shared formula in C3:AI3
# recalculates off of it automatically.
#
# The rest of the commit is just the file having been re-saved by a
# newer Excel build, which also nudged the saved cursor/selection and
# active sheet/tab - reproduce those view-state tweaks too.

$wb  = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsIEA   = $wb.Worksheets.Item("IEA Data")
$wsHPEbP = $wb.Worksheets.Item("HPEbP")

# --- the actual data edit -------------------------------------------------
$wsHPEbP.Range("B3").Formula = "=118/(162+2)"

# --- saved view state (selection per sheet + active sheet/tab) -----------
$wsAbout.Range("B14").Select()
$wsIEA.Range("D7:F7").Select()
$wsHPEbP.Range("C3").Select()

# HPEbP is the sheet that is active/selected when the file is saved
$wsHPEbP.Activate()
